$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: Mirko, 10/01/2017 (serial 42745), "calcolatrice" activity, 4 hours (0.16666...)
$ws.Cells.Item(24, 1).Value = 42745
$ws.Cells.Item(24, 2).Value = "Mirko"
$ws.Cells.Item(24, 3).Value = "Implementatata bozza grafica calcolatrice con 6 bottoni, testo e icone"
$ws.Cells.Item(24, 4).Value = 0.16666666666666666

# Highlight the grand-total cell (I2) with an underline + elapsed-time format,
# mirroring the new style added to the workbook (new font + cellXf).
$ws.Range("I2").Font.Underline = $true
$ws.Range("I2").NumberFormat = "[h]:mm:ss"

# Recalculate so cached formula results (F2/G2/I2 totals) reflect the new row.
$excel.Calculate()
